$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing AgTests (F) / AgPosit (G) values for rows 492-525 ---
$ws.Range("F492").Value = 14292

$ws.Range("F494").Value = 6712
$ws.Range("F495").Value = 10426
$ws.Range("F496").Value = 8251
$ws.Range("F497").Value = 7740
$ws.Range("F498").Value = 9158
$ws.Range("F499").Value = 11418
$ws.Range("F500").Value = 7726
$ws.Range("F501").Value = 5738
$ws.Range("F502").Value = 10539
$ws.Range("F503").Value = 7460
$ws.Range("F504").Value = 7487
$ws.Range("F505").Value = 8528
$ws.Range("F506").Value = 10843
$ws.Range("F507").Value = 7219
$ws.Range("F508").Value = 5673
$ws.Range("F509").Value = 9559
$ws.Range("F510").Value = 7846
$ws.Range("F511").Value = 6791
$ws.Range("F512").Value = 8484
$ws.Range("F513").Value = 10353
$ws.Range("F514").Value = 6946

$ws.Range("F515").Value = 5006
$ws.Range("G515").Value = 16

$ws.Range("F516").Value = 9297
$ws.Range("F517").Value = 6712
$ws.Range("F518").Value = 7013

$ws.Range("F519").Value = 7816
$ws.Range("G519").Value = 21

$ws.Range("F520").Value = 10041
$ws.Range("F521").Value = 6445
$ws.Range("F522").Value = 4969
$ws.Range("F523").Value = 9724

$ws.Range("F524").Value = 7557
$ws.Range("G524").Value = 29

$ws.Range("F525").Value = 7247
$ws.Range("G525").Value = 21

# --- Append new rows 526-529 ---
$ws.Range("A526").Value = 44420
$ws.Range("B526").Value = 393361
$ws.Range("C526").Value = 6098
$ws.Range("D526").Value = 59
$ws.Range("E526").Value = 12544
$ws.Range("F526").Value = 8232
$ws.Range("G526").Value = 24

$ws.Range("A527").Value = 44421
$ws.Range("B527").Value = 393455
$ws.Range("C527").Value = 9202
$ws.Range("D527").Value = 94
$ws.Range("E527").Value = 12544
$ws.Range("F527").Value = 10049
$ws.Range("G527").Value = 31

$ws.Range("A528").Value = 44422
$ws.Range("B528").Value = 393529
$ws.Range("C528").Value = 4692
$ws.Range("D528").Value = 74
$ws.Range("E528").Value = 12544
$ws.Range("F528").Value = 6402
$ws.Range("G528").Value = 16

$ws.Range("A529").Value = 44423
$ws.Range("B529").Value = 393536
$ws.Range("C529").Value = 1309
$ws.Range("D529").Value = 7
$ws.Range("E529").Value = 12544
$ws.Range("F529").Value = 3650
$ws.Range("G529").Value = 18
